# Update crypto price/volume data per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.788.86"
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("D3").Value = "1.911.91"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'250.72"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'0.703"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'46.84"
$ws.Range("E8").Value = "  +8.08%  "
$ws.Range("D9").Value = "'0.373"
$ws.Range("E9").Value = "  +4.29%  "
$ws.Range("E10").Value = "  +9.51%  "
$ws.Range("D11").Value = "'0.0763"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "'14.62"
$ws.Range("E13").Value = "  +7.99%  "
$ws.Range("D14").Value = "'0.817"
$ws.Range("E14").Value = "  +5.44%  "
$ws.Range("D15").Value = "2.196.82"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "'5.14"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("D17").Value = "1.907.81"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "37.144.90"
$ws.Range("E18").Value = "  +4.77%  "
$ws.Range("D19").Value = "'74.79"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "0.0₃0861"
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("E21").Value = "  +6.46%  "
$ws.Range("D22").Value = "'251.55"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "'5.19"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "'2.62"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("D27").Value = "'167.82"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "'8.83"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "'18.72"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("E31").Value = "  +6.77%  "
$ws.Range("D32").Value = "'0.0620"
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("D33").Value = "'0.0911"
$ws.Range("E33").Value = "  +23.50%  "
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("D35").Value = "'1.90"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  +4.97%  "
$ws.Range("D38").Value = "'18.66"
$ws.Range("E38").Value = "  +55.34%  "
$ws.Range("D39").Value = "'0.875"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("D41").Value = "'105.69"
$ws.Range("E41").Value = "  +8.41%  "
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("D44").Value = "'2.86"
$ws.Range("E44").Value = "  +19.43%  "
$ws.Range("D45").Value = "'1.10"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("D46").Value = "1.349.47"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'0.0812"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").Value = "'6.48"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("D51").Value = "'43.30"
$ws.Range("E51").Value = "  +1.18%  "
